$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    "C2" = 4.026394718256862
    "D2" = 9.960504504094272
    "E2" = 13.87345592677048
    "F2" = 28.5022060582568
    "G2" = 27.0497316649349
    "H2" = 13.99239377049987
    "I2" = 19.62191857968765
    "J2" = 9.905874282657706
    "K2" = 13.85906985517805
    "M2" = 16.63564333578771
    "O2" = 21.04636493116709
    "C3" = 3.87836050223641
    "D3" = 9.938361622787882
    "E3" = 13.90408817442267
    "F3" = 28.6382390188574
    "G3" = 27.23842511706677
    "H3" = 14.06051134251524
    "I3" = 19.72102959989156
    "J3" = 9.940445628200685
    "K3" = 13.23564009844841
    "M3" = 16.37230878384486
    "O3" = 21.1712078546357
    "C4" = 3.783749376752831
    "D4" = 9.926154998038996
    "E4" = 13.92541263193496
    "F4" = 28.72988042724842
    "G4" = 27.36538168959237
    "H4" = 14.10498103141494
    "I4" = 19.7872221408544
    "J4" = 9.962947605459032
    "K4" = 12.83673264260108
    "M4" = 16.20949270048343
    "O4" = 21.25333432047183
    "C5" = 3.744294724026761
    "D5" = 9.921533513642657
    "E5" = 13.93473480297431
    "F5" = 28.76925900776051
    "G5" = 27.41989000523393
    "H5" = 14.12376817848645
    "I5" = 19.81553453707211
    "J5" = 9.972438480304193
    "K5" = 12.67029968312607
    "M5" = 16.14292980457407
    "O5" = 21.28817512560725
    "C6" = 3.737690065260456
    "D6" = 9.920787525828
    "E6" = 13.93632091677255
    "F6" = 28.77592042989533
    "G6" = 27.42910798770389
    "H6" = 14.12692796484539
    "I6" = 19.82031651115187
    "J6" = 9.974033844078269
    "K6" = 12.64243532533211
    "M6" = 16.13186621502252
    "O6" = 21.2940433059397
    "C7" = 3.783220871678119
    "D7" = 9.926091238124251
    "E7" = 13.92553579452218
    "F7" = 28.7304032744963
    "G7" = 27.36610560623233
    "H7" = 14.10523170700979
    "I7" = 19.78759855775333
    "J7" = 9.963074301754119
    "K7" = 12.8345035035107
    "M7" = 16.20859578383073
    "O7" = 21.25379863732135
    "C8" = 3.976143763674748
    "D8" = 9.952583477685787
    "E8" = 13.88349547176233
    "F8" = 28.54742142858409
    "G8" = 27.11247883856782
    "H8" = 14.01533182559487
    "I8" = 19.65498183983946
    "J8" = 9.917530162461881
    "K8" = 13.6475556530895
    "M8" = 16.54512092644807
    "O8" = 21.08827347102956
    "C9" = 4.323640392585195
    "D9" = 10.01539420734227
    "E9" = 13.82103536697227
    "F9" = 28.25330548266389
    "G9" = 26.7039861489207
    "H9" = 13.86001675586226
    "I9" = 19.4374474492059
    "J9" = 9.838311357682443
    "K9" = 15.10749607952762
    "M9" = 17.19311817301924
    "O9" = 20.80721227081638
    "C10" = 4.558641431125698
    "D10" = 10.06792899429836
    "E10" = 13.78734565809013
    "F10" = 28.07708297151467
    "G10" = 26.45914580997697
    "H10" = 13.75868120192179
    "I10" = 19.30378388430997
    "J10" = 9.786229218723969
    "K10" = 16.0909322838979
    "M10" = 17.65806886997544
    "O10" = 20.62741115557801
    "C11" = 4.660887784320527
    "D11" = 10.09316104925264
    "E11" = 13.77467067168753
    "F11" = 28.00565707934643
    "G11" = 26.35999883489184
    "H11" = 13.71535240102514
    "I11" = 19.24870661282336
    "J11" = 9.763857676511922
    "K11" = 16.51781431951597
    "M11" = 17.86639576976022
    "O11" = 20.55144650603362
    "C12" = 4.698918865248034
    "D12" = 10.10290256568552
    "E12" = 13.77025215224811
    "F12" = 27.9798734522086
    "G12" = 26.32423250848062
    "H12" = 13.69934308323065
    "I12" = 19.22867791239708
    "J12" = 9.755575604406515
    "K12" = 16.67644242093916
    "M12" = 17.9447650805198
    "O12" = 20.52352187432347
    "C13" = 4.690759019255138
    "D13" = 10.10079633648356
    "E13" = 13.7711868045632
    "F13" = 27.98537009679412
    "G13" = 26.33185602031817
    "H13" = 13.7027732564908
    "I13" = 19.23295456324655
    "J13" = 9.757350872932623
    "K13" = 16.64241454580305
    "M13" = 17.9279109122277
    "O13" = 20.52949845930735
    "C14" = 4.664030493280618
    "D14" = 10.09395877821577
    "E14" = 13.7742995171607
    "F14" = 28.00351047565147
    "G14" = 26.35702056338542
    "H14" = 13.71402732016667
    "I14" = 19.24704222322836
    "J14" = 9.763172508914824
    "K14" = 16.53092572414673
    "M14" = 17.8728540181248
    "O14" = 20.54913224383809
    "C15" = 4.647568497090651
    "D15" = 10.08979473313045
    "E15" = 13.77625579011336
    "F15" = 28.01478677713757
    "G15" = 26.37266675835776
    "H15" = 13.72097264204662
    "I15" = 19.25577926235618
    "J15" = 9.766763101880247
    "K15" = 16.46223988514427
    "M15" = 17.83906062566082
    "O15" = 20.56126820160241
    "C16" = 4.551864108206641
    "D16" = 10.06630642678235
    "E16" = 13.78822734371076
    "F16" = 28.08192728519702
    "G16" = 26.4658731134457
    "H16" = 13.76156857687893
    "I16" = 19.30749892640901
    "J16" = 9.787717797083673
    "K16" = 16.06261540300494
    "M16" = 17.6443848253616
    "O16" = 20.63249319156843
    "C17" = 4.491946386076321
    "D17" = 10.05223514783111
    "E17" = 13.79625047392812
    "F17" = 28.12535916854366
    "G17" = 26.52619909484341
    "H17" = 13.78718223105373
    "I17" = 19.34069717697909
    "J17" = 9.800910848787396
    "K17" = 15.81215233809672
    "M17" = 17.52409759454874
    "O17" = 20.67768251692899
    "C18" = 4.457045863574198
    "D18" = 10.0442675755671
    "E18" = 13.80111465883382
    "F18" = 28.15116213931396
    "G18" = 26.56204719536378
    "H18" = 13.80217517169607
    "I18" = 19.36033072407576
    "J18" = 9.808623494053684
    "K18" = 15.66616783689462
    "M18" = 17.4546156369954
    "O18" = 20.70422254488533
    "C19" = 4.445154600256931
    "D19" = 10.04159166013646
    "E19" = 13.80280443515863
    "F19" = 28.16003955566267
    "G19" = 26.57438169421567
    "H19" = 13.8072962858707
    "I19" = 19.36707071078524
    "J19" = 9.811256232846922
    "K19" = 15.61641197126783
    "M19" = 17.43104135324318
    "O19" = 20.71330261199378
    "C20" = 4.498370139354642
    "D20" = 10.05372007022967
    "E20" = 13.79537057566758
    "F20" = 28.1206506333021
    "G20" = 26.51965810845949
    "H20" = 13.7844286372007
    "I20" = 19.33710737205559
    "J20" = 9.799493558963098
    "K20" = 15.83901436343578
    "M20" = 17.53693346141497
    "O20" = 20.67281526157072
    "C21" = 4.671900093810169
    "D21" = 10.09596210997595
    "E21" = 13.77337489174153
    "F21" = 27.99814785080874
    "G21" = 26.34958070532662
    "H21" = 13.71071091620277
    "I21" = 19.24288183360519
    "J21" = 9.761457412204972
    "K21" = 16.56375527927373
    "M21" = 17.88904014354777
    "O21" = 20.54334245809065
    "C22" = 4.781298796410926
    "D22" = 10.12465535158808
    "E22" = 13.7612214545339
    "F22" = 27.9254548420887
    "G22" = 26.24880050591438
    "H22" = 13.66485415498778
    "I22" = 19.18612741544345
    "J22" = 9.73770323085393
    "K22" = 17.01976914907187
    "M22" = 18.11610708860263
    "O22" = 20.46363127756784
    "C23" = 4.723283073529048
    "D23" = 10.10924363394707
    "E23" = 13.76750466404989
    "F23" = 27.96357581715536
    "G23" = 26.30163308183127
    "H23" = 13.68911624901443
    "I23" = 19.2159751874482
    "J23" = 9.75028033806444
    "K23" = 16.77802261658014
    "M23" = 17.99521647618406
    "O23" = 20.50572448411368
    "C24" = 4.495467370212287
    "D24" = 10.05304835581891
    "E24" = 13.79576759403113
    "F24" = 28.12277676663383
    "G24" = 26.5226116582064
    "H24" = 13.78567270432811
    "I24" = 19.33872861797373
    "J24" = 9.800133917620984
    "K24" = 15.82687623408976
    "M24" = 17.53113138049264
    "O24" = 20.67501400334383
    "C25" = 4.233103060062366
    "D25" = 9.997262680743086
    "E25" = 13.8357910805628
    "F25" = 28.32590049560651
    "G25" = 26.8048660101242
    "H25" = 13.89978934824318
    "I25" = 19.49171927509168
    "J25" = 9.858665067304919
    "K25" = 14.72776050922688
    "M25" = 17.01950163286428
    "O25" = 20.8785704997411
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
